$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Fix the CasesTab query (cell B2): remove the trailing Cohort column from the
# RETURN clause since the `cohort` node is no longer part of the result set.
$newQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)" + [char]10 + `
"MATCH (c)<--(diag:diagnosis)" + [char]10 + `
"OPTIONAL MATCH (samp:sample)-->(c)" + [char]10 + `
"OPTIONAL MATCH (co:cohort)<-[*]-(c)" + [char]10 + `
"WITH DISTINCT c, s, demo, diag, co" + [char]10 + `
"WHERE diag.primary_disease_site IN ['Bone']" + [char]10 + `
"RETURN  coalesce(c.case_id, '') AS ``Case ID`` ," + [char]10 + `
"        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ," + [char]10 + `
"        coalesce(s.clinical_study_type, '') AS  ``Study Type``," + [char]10 + `
"        coalesce(demo.breed, '') AS Breed ," + [char]10 + `
"        coalesce(diag.disease_term, '') AS Diagnosis ," + [char]10 + `
"        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ," + [char]10 + `
"        coalesce(demo.patient_age_at_enrollment, '') AS Age ," + [char]10 + `
"        coalesce(demo.sex, '') AS Sex ," + [char]10 + `
"        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``," + [char]10 + `
"        coalesce(demo.weight, '') AS ``Weight (kg)``," + [char]10 + `
"        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newQuery

# The RETURN clause got shorter, so the word-wrapped rows re-flow to a
# shorter auto-fit height.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# Update the active selection to match the saved workbook state (B2 selected,
# scrolled back to the top of the sheet).
$ws.Range("B2").Select()
